# Apply the weekly report refresh edit:
#  - Update the "Report Generated On" timestamp
#  - Zero out the billed amount / pricing total cells (no violations this run)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# Update report generation timestamp
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"

# Zero out Total Billed Amount
$ws.Range("C8").Value = 0

# Zero out per-line-item pricing and section totals in the H column
$zeroRows = @(16, 17, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36)
foreach ($r in $zeroRows) {
    $ws.Range("H$r").Value = 0
}
